$d = $word.ActiveDocument

# Start from the end of the last existing paragraph in the document
$tail = $d.Paragraphs.Last.Range
$tail.Collapse(0)

# --- New list paragraph 1 ---
$tail.InsertParagraphAfter()
# Re-fetch the range of the newly created (now last) paragraph; the
# previous Range object does not track past the inserted paragraph mark.
$tail = $d.Paragraphs.Last.Range
$tail.Collapse(0)
$tail.InsertAfter("Dependencies are of two types: regular dependencies where the code depends on the dependency and the dev-dependencies that help the developers with tools and helpers for the developers(like webpack, nodemon etc.)")
$tail.Collapse(0)
$tail.InsertAfter(". Nodemon is good for listening for the changes and restarting the server for us instead of using “node index.js”.")
$tail.Collapse(0)

# --- New list paragraph 2 ---
$tail.InsertParagraphAfter()
# Re-fetch the range of the newly created (now last) paragraph; the
# previous Range object does not track past the inserted paragraph mark.
$tail = $d.Paragraphs.Last.Range
$tail.Collapse(0)
$tail.InsertAfter("If we want to use a local script in the terminal we need to declare them into the scripts. We can only use global dependencies directly in the terminal")
$tail.Collapse(0)

# --- New list paragraph 3 ---
$tail.InsertParagraphAfter()
# Re-fetch the range of the newly created (now last) paragraph; the
# previous Range object does not track past the inserted paragraph mark.
$tail = $d.Paragraphs.Last.Range
$tail.Collapse(0)
$tail.InsertAfter("What is a “slug”? It is")
$tail.Collapse(0)
$tail.InsertAfter(" a")
$tail.Collapse(0)
$tail.InsertAfter(" part of a URL that identifies a particular page on a website in an easy-to-read form.")
$tail.Collapse(0)
$tail.InsertAfter(" ")
$tail.Collapse(0)
$tail.InsertAfter("A slug is the part of the URL that explains the page’s content. ")
$tail.Collapse(0)
$tail.InsertAfter(" (slugify used)")
$tail.Collapse(0)

# --- New list paragraph 4 ---
$tail.InsertParagraphAfter()
# Re-fetch the range of the newly created (now last) paragraph; the
# previous Range object does not track past the inserted paragraph mark.
$tail = $d.Paragraphs.Last.Range
$tail.Collapse(0)
$tail.InsertAfter("Code for dependencies versions. Ex: ")
$tail.Collapse(0)
$tail.InsertAfter("*")
$tail.Collapse(0)
$tail.InsertAfter("^")
$tail.Collapse(0)
$tail.InsertAfter("~")
$tail.Collapse(0)
$tail.InsertAfter("1.18.10 (1 is the major version, 18 is the minor version, and the 10 is the patch version)")
$tail.Collapse(0)
$tail.InsertAfter(". The “~” stands for only patch releases(which is safer")
$tail.Collapse(0)
$tail.InsertAfter(")")
$tail.Collapse(0)
$tail.InsertAfter(", the “^” is for all the patch and minor releases")
$tail.Collapse(0)
$tail.InsertAfter(". There is also an option to update to all version using “*”")
$tail.Collapse(0)

# --- New list paragraph 5 ---
$tail.InsertParagraphAfter()
# Re-fetch the range of the newly created (now last) paragraph; the
# previous Range object does not track past the inserted paragraph mark.
$tail = $d.Paragraphs.Last.Range
$tail.Collapse(0)
$tail.InsertAfter("Util commands : npm outdated")
$tail.Collapse(0)
$tail.InsertAfter(", npm update(that doesn’t work")
$tail.Collapse(0)
$tail.InsertAfter(" updating the package.json, but package-lock is working ")
$tail.Collapse(0)
$tail.InsertAfter("in npm 6+), npm I package@version, npm update <packageName>")
$tail.Collapse(0)
$tail.InsertAfter(", npm -rm -r node_modules, npm uninstall <packageName>")
$tail.Collapse(0)
$tail.InsertAfter(", npm install")
$tail.Collapse(0)

Write-Output ("Inserted paragraphs. Total paragraphs now: " + $d.Paragraphs.Count)
